# Server to client broadcasts for day
# Add "Vote Start" and "Vote Received" message rows to the Vote section of the
# message API sheet, and clarify the Result row's comment (map -> array).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for "Vote Start" right above the existing "Vote" row ---
$ws.Rows.Item(35).Insert()
$ws.Range("C35").Value = "Vote Start"
$ws.Range("D35").Value = "voteStart"

# --- Insert a new row for "Vote Received" right below the "Vote" row       ---
# (the "Vote" row has shifted from 35 -> 36 because of the insert above, and
# the old "Result" row has shifted from 36 -> 37)
$ws.Rows.Item(37).Insert()
$ws.Range("C37").Value = "Vote Received"
$ws.Range("D37").Value = "voteReceived"
$ws.Range("E37").Value = "playerID"
$ws.Range("F37").Value = "Player ID that sent the vote; broadcast to all"

# --- Update the Result row comment: map -> array                          ---
$ws.Range("F38").Value = "Votes = array from id to votes, playerRoles = array from id to roles, winTeam = id"

# --- Restore the saved selection/scroll position                          ---
$ws.Range("F39").Select()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
